# Update "想去人数" (column F) counts for several events that appear
# both on the "展览" (Exhibition) sheet and the combined "全部类型"
# (All Types) sheet, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 2691
$wsExhibit.Range("F10").Value = 1566
$wsExhibit.Range("F17").Value = 6
$wsExhibit.Range("F24").Value = 1732
$wsExhibit.Range("F27").Value = 66
$wsExhibit.Range("F31").Value = 442

# Sheet "全部类型": same events, shifted one row down because it also
# contains the single "演出" entry.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 2691
$wsAll.Range("F11").Value = 1566
$wsAll.Range("F18").Value = 6
$wsAll.Range("F25").Value = 1732
$wsAll.Range("F28").Value = 66
$wsAll.Range("F32").Value = 442
